$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Todo")
$ws.Activate()

# Status column (D) - mark first few tasks as "Completeish"
$ws.Range("D2").Value = "Completeish"
$ws.Range("D3").Value = "Completeish"
$ws.Range("D4").Value = "Completeish"

# Owner column (C) - written in the order needed to reproduce the
# shared-string table ordering from the original authoring session.
$ws.Range("C5").Value = "Jameson Riley"
$ws.Range("C6").Value = "Jameson Riley"
$ws.Range("C14").Value = "Jameson Riley"
$ws.Range("C9").Value = "Naseem"
$ws.Range("C10").Value = "Steven"
$ws.Range("C11").Value = "Steven"
$ws.Range("C12").Value = "Steven"
$ws.Range("C13").Value = "Steven"
$ws.Range("C7").Value = "Andrew"
$ws.Range("C8").Value = "Naseem/Andrew"

# Column widths (B=78, C~13.57 best-fit, D=12)
$ws.Columns.Item(2).ColumnWidth = 77.15
$ws.Columns.Item(3).ColumnWidth = 12.6
$ws.Columns.Item(4).ColumnWidth = 11.15

# View / selection changes
$ws.Range("B11").Select()
$excel.ActiveWindow.Zoom = 190
